$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1661084259327765
$ws.Range("D2").Value = 0.005612785346094995
$ws.Range("E2").Value = 0.1074619746187437
$ws.Range("F2").Value = 0.3746353073733673
$ws.Range("G2").Value = 0.2608699164456141
$ws.Range("H2").Value = 0.3328812873599105
$ws.Range("I2").Value = 0.3463768926081165
$ws.Range("M2").Value = 0.5398403867358468
$ws.Range("N2").Value = 1.251098144649433
$ws.Range("O2").Value = 1.103858612678096

$ws.Range("B3").Value = 0.1472636846431499
$ws.Range("D3").Value = 0.004994341764724908
$ws.Range("E3").Value = 0.1092998160422169
$ws.Range("F3").Value = 0.3525363198711773
$ws.Range("G3").Value = 0.2385975703545427
$ws.Range("H3").Value = 0.3259847953969484
$ws.Range("I3").Value = 0.3463474110070166
$ws.Range("M3").Value = 0.4722893597765676
$ws.Range("N3").Value = 1.210295828351775
$ws.Range("O3").Value = 1.04189678590032

$ws.Range("B4").Value = 0.1356685332018657
$ws.Range("D4").Value = 0.004615247115150822
$ws.Range("E4").Value = 0.1106025869431377
$ws.Range("F4").Value = 0.3392116276521406
$ws.Range("G4").Value = 0.2250785283969776
$ws.Range("H4").Value = 0.3219546496304275
$ws.Range("I4").Value = 0.346545782202007
$ws.Range("M4").Value = 0.4307244345333032
$ws.Range("N4").Value = 1.185635250531561
$ws.Range("O4").Value = 1.004599322642918

$ws.Range("B5").Value = 0.1309377210038747
$ws.Range("D5").Value = 0.004460933676792678
$ws.Range("E5").Value = 0.1111769599516208
$ws.Range("F5").Value = 0.3338429597921007
$ws.Range("G5").Value = 0.2196084561509934
$ws.Range("H5").Value = 0.3203636907477687
$ws.Range("I5").Value = 0.3466812631422691
$ws.Range("M5").Value = 0.4137649062256656
$ws.Range("N5").Value = 1.175686032843274
$ws.Range("O5").Value = 0.9895876508836636

$ws.Range("B6").Value = 0.1301518424028387
$ws.Range("D6").Value = 0.004435320703166212
$ws.Range("E6").Value = 0.1112749501050736
$ws.Range("F6").Value = 0.3329551917873985
$ws.Range("G6").Value = 0.2187025059342744
$ws.Range("H6").Value = 0.3201026145344628
$ws.Range("I6").Value = 0.3467070661298024
$ws.Range("M6").Value = 0.4109475133614922
$ws.Range("N6").Value = 1.174040068323421
$ws.Range("O6").Value = 0.9871062726617197

$ws.Range("B7").Value = 0.1356047542801946
$ws.Range("D7").Value = 0.00461316528325284
$ws.Range("E7").Value = 0.1106101575437499
$ws.Range("F7").Value = 0.3391389761076908
$ws.Range("G7").Value = 0.2250045993183534
$ws.Range("H7").Value = 0.3219329854854607
$ws.Range("I7").Value = 0.346547387818525
$ws.Range("M7").Value = 0.4304957984247295
$ws.Range("N7").Value = 1.18550066427126
$ws.Range("O7").Value = 1.004396112196787

$ws.Range("B8").Value = 0.1596160835493521
$ws.Range("D8").Value = 0.00539942163900875
$ws.Range("E8").Value = 0.1080592679174188
$ws.Range("F8").Value = 0.3669648016594351
$ws.Range("G8").Value = 0.2531578219523141
$ws.Range("H8").Value = 0.3304609523480053
$ws.Range("I8").Value = 0.3463219221861387
$ws.Range("M8").Value = 0.5165675376332928
$ws.Range("N8").Value = 1.236949010496573
$ws.Range("O8").Value = 1.082338692579725

$ws.Range("B9").Value = 0.2064914003842091
$ws.Range("D9").Value = 0.006945842567709803
$ws.Range("E9").Value = 0.1044560952958236
$ws.Range("F9").Value = 0.4234784477749542
$ws.Range("G9").Value = 0.3096202400279822
$ws.Range("H9").Value = 0.3488078672512955
$ws.Range("I9").Value = 0.3475887679484373
$ws.Range("M9").Value = 0.6846295097790005
$ws.Range("N9").Value = 1.340889488402212
$ws.Range("O9").Value = 1.241148948972636

$ws.Range("B10").Value = 0.2407823869148444
$ws.Range("D10").Value = 0.008084302042611569
$ws.Range("E10").Value = 0.1026837273292429
$ws.Range("F10").Value = 0.4662044215769754
$ws.Range("G10").Value = 0.3518914261744612
$ws.Range("H10").Value = 0.3632823492797286
$ws.Range("I10").Value = 0.3495505719500969
$ws.Range("M10").Value = 0.8076433401380854
$ws.Range("N10").Value = 1.419040578161599
$ws.Range("O10").Value = 1.361525135193943

$ws.Range("B11").Value = 0.2563460674260511
$ws.Range("D11").Value = 0.008602620141545003
$ws.Range("E11").Value = 0.1020721353747511
$ws.Range("F11").Value = 0.4859075045388863
$ws.Range("G11").Value = 0.3712981816973979
$ws.Range("H11").Value = 0.3700844682704343
$ws.Range("I11").Value = 0.3506647198888473
$ws.Range("M11").Value = 0.8635019757502107
$ws.Range("N11").Value = 1.454966479155587
$ws.Range("O11").Value = 1.417104467211857

$ws.Range("B12").Value = 0.2622341251766898
$ws.Range("D12").Value = 0.008798944433408451
$ws.Range("E12").Value = 0.1018689220452096
$ws.Range("F12").Value = 0.4934071510788414
$ws.Range("G12").Value = 0.3786728547891585
$ws.Range("H12").Value = 0.3726916047249631
$ws.Range("I12").Value = 0.3511183062409131
$ws.Range("M12").Value = 0.8846391120588493
$ws.Range("N12").Value = 1.468623137977346
$ws.Range("O12").Value = 1.438269610491318

$ws.Range("B13").Value = 0.2609662818065601
$ws.Range("D13").Value = 0.008756660535020444
$ws.Range("E13").Value = 0.1019114190903103
$ws.Range("F13").Value = 0.4917902526312048
$ws.Range("G13").Value = 0.3770834375360437
$ws.Range("H13").Value = 0.3721287177191073
$ws.Range("I13").Value = 0.3510192123658271
$ws.Range("M13").Value = 0.8800875454198405
$ws.Range("N13").Value = 1.465679630869687
$ws.Range("O13").Value = 1.433706040295874

$ws.Range("B14").Value = 0.2568305958083954
$ws.Range("D14").Value = 0.008618770961767552
$ws.Range("E14").Value = 0.10205484580883
$ws.Range("F14").Value = 0.486523732200169
$ws.Range("G14").Value = 0.3719043828823203
$ws.Range("H14").Value = 0.3702983310493551
$ws.Range("I14").Value = 0.350701402878812
$ws.Range("M14").Value = 0.8652412521942097
$ws.Range("N14").Value = 1.456088981442775
$ws.Range("O14").Value = 1.418843358821618

$ws.Range("B15").Value = 0.2542966293602547
$ws.Range("D15").Value = 0.008534315545453808
$ws.Range("E15").Value = 0.1021464067595801
$ws.Range("F15").Value = 0.4833028567377369
$ws.Range("G15").Value = 0.3687354239147567
$ws.Range("H15").Value = 0.3691812462652564
$ws.Range("I15").Value = 0.3505108555245258
$ws.Range("M15").Value = 0.8561454506516952
$ws.Range("N15").Value = 1.450221196524694
$ws.Range("O15").Value = 1.409754981530398

$ws.Range("B16").Value = 0.2397645099099179
$ws.Range("D16").Value = 0.008050436224955604
$ws.Range("E16").Value = 0.1027276468475939
$ws.Range("F16").Value = 0.4649221605921383
$ws.Range("G16").Value = 0.3506267397123963
$ws.Range("H16").Value = 0.3628421961578994
$ws.Range("I16").Value = 0.3494822041198944
$ws.Range("M16").Value = 0.8039907468769769
$ws.Range("N16").Value = 1.416700139513921
$ws.Range("O16").Value = 1.357909438847571

$ws.Range("B17").Value = 0.2308401070648642
$ws.Range("D17").Value = 0.007753692137388413
$ws.Range("E17").Value = 0.1031343551526405
$ws.Range("F17").Value = 0.4537146441983424
$ws.Range("G17").Value = 0.3395632573981402
$ws.Range("H17").Value = 0.3590091495388776
$ws.Range("I17").Value = 0.3489078035303805
$ws.Range("M17").Value = 0.7719690954080534
$ws.Range("N17").Value = 1.3962308782489
$ws.Range("O17").Value = 1.326314224011327

$ws.Range("B18").Value = 0.2257037168302531
$ws.Range("D18").Value = 0.007583053770687798
$ws.Range("E18").Value = 0.1033865825994837
$ws.Range("F18").Value = 0.4472934830934321
$ws.Range("G18").Value = 0.333216516122846
$ws.Range("H18").Value = 0.3568249626434294
$ws.Range("I18").Value = 0.3485983040424045
$ws.Range("M18").Value = 0.7535416317556383
$ws.Range("N18").Value = 1.384492850305435
$ws.Range("O18").Value = 1.3082186011749

$ws.Range("B19").Value = 0.2239640692837384
$ws.Range("D19").Value = 0.007525285984804952
$ws.Range("E19").Value = 0.1034751135740812
$ws.Range("F19").Value = 0.4451236957152673
$ws.Range("G19").Value = 0.3310704759802121
$ws.Range("H19").Value = 0.3560889510529393
$ws.Range("I19").Value = 0.3484971051891819
$ws.Range("M19").Value = 0.7473008093512874
$ws.Range("N19").Value = 1.380524674956973
$ws.Range("O19").Value = 1.30210495948748

$ws.Range("B20").Value = 0.2317904709482264
$ws.Range("D20").Value = 0.007785276904343164
$ws.Range("E20").Value = 0.1030891635963407
$ws.Range("F20").Value = 0.4549051040104644
$ws.Range("G20").Value = 0.3407392557370343
$ws.Range("H20").Value = 0.3594150640562361
$ws.Range("I20").Value = 0.3489667898463011
$ws.Range("M20").Value = 0.7753788393913794
$ws.Range("N20").Value = 1.398406218859037
$ws.Range("O20").Value = 1.32966960372076

$ws.Range("B21").Value = 0.2580455019559906
$ws.Range("D21").Value = 0.008659271259254808
$ws.Range("E21").Value = 0.1020119444744836
$ws.Range("F21").Value = 0.4880695904105181
$ws.Range("G21").Value = 0.3734248957898103
$ws.Range("H21").Value = 0.3708351096815505
$ws.Range("I21").Value = 0.3507938929506693
$ws.Range("M21").Value = 0.869602390414471
$ws.Range("N21").Value = 1.458904580659038
$ws.Range("O21").Value = 1.423205669011054

$ws.Range("B22").Value = 0.2751719974716309
$ws.Range("D22").Value = 0.009230753626781052
$ws.Range("E22").Value = 0.1014735408759151
$ws.Range("F22").Value = 0.5099690716487117
$ws.Range("G22").Value = 0.394937091137507
$ws.Range("H22").Value = 0.3784813502030033
$ws.Range("I22").Value = 0.3521725901145345
$ws.Range("M22").Value = 0.9310932837195907
$ws.Range("N22").Value = 1.498747884015131
$ws.Range("O22").Value = 1.48502762612577

$ws.Range("B23").Value = 0.2660344102659451
$ws.Range("D23").Value = 0.008925721714049928
$ws.Range("E23").Value = 0.1017456133722199
$ws.Range("F23").Value = 0.4982603094591411
$ws.Range("G23").Value = 0.3834417987347365
$ws.Range("H23").Value = 0.374383691544466
$ws.Range("I23").Value = 0.3514199299106977
$ws.Range("M23").Value = 0.8982829013134221
$ws.Range("N23").Value = 1.477455448994789
$ws.Range("O23").Value = 1.451968684885117

$ws.Range("B24").Value = 0.2313608289441333
$ws.Range("D24").Value = 0.007770997543271818
$ws.Range("E24").Value = 0.1031095373970317
$ws.Range("F24").Value = 0.4543668280221738
$ws.Range("G24").Value = 0.3402075439158239
$ws.Range("H24").Value = 0.3592314894154072
$ws.Range("I24").Value = 0.3489400575431318
$ws.Range("M24").Value = 0.7738373494950537
$ws.Range("N24").Value = 1.397422653730075
$ws.Range("O24").Value = 1.328152422018917

$ws.Range("B25").Value = 0.1938350087837364
$ws.Range("D25").Value = 0.006527055015524752
$ws.Range("E25").Value = 0.1052787998005797
$ws.Range("F25").Value = 0.4079796454635698
$ws.Range("G25").Value = 0.2942089477174648
$ws.Range("H25").Value = 0.3436701850168191
$ws.Range("I25").Value = 0.3470642160299064
$ws.Range("M25").Value = 0.6392438986159306
$ws.Range("N25").Value = 1.312451925063584
$ws.Range("O25").Value = 1.197541566728916

